$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 14.70619329725741
$ws.Range("C2").Value = 14.4909633698582
$ws.Range("D2").Value = 4.477272852686781
$ws.Range("F2").Value = 21.33687627677802
$ws.Range("G2").Value = 23.67168218289414
$ws.Range("H2").Value = 13.03675337168901
$ws.Range("L2").Value = 10.96841794683468
$ws.Range("M2").Value = 14.74867058625313
$ws.Range("O2").Value = 19.15831243127248
# Row 3
$ws.Range("B3").Value = 14.07629987912296
$ws.Range("C3").Value = 14.37572802513451
$ws.Range("D3").Value = 4.402320590760076
$ws.Range("F3").Value = 21.39530775670628
$ws.Range("G3").Value = 23.76620978859924
$ws.Range("H3").Value = 13.09548705294023
$ws.Range("L3").Value = 10.99141817284869
$ws.Range("M3").Value = 14.61902162110279
$ws.Range("O3").Value = 19.25513350257224
# Row 4
$ws.Range("B4").Value = 13.6750328188273
$ws.Range("C4").Value = 14.305666470245
$ws.Range("D4").Value = 4.355212562583664
$ws.Range("F4").Value = 21.43867056577712
$ws.Range("G4").Value = 23.83547918707503
$ws.Range("H4").Value = 13.13418788079078
$ws.Range("L4").Value = 11.00732654062118
$ws.Range("M4").Value = 14.54048877328303
$ws.Range("O4").Value = 19.32006879799615
# Row 5
$ws.Range("B5").Value = 13.50807138925405
$ws.Range("C5").Value = 14.27731210105228
$ws.Range("D5").Value = 4.335756732134951
$ws.Range("F5").Value = 21.45821541247282
$ws.Range("G5").Value = 23.86650950566243
$ws.Range("H5").Value = 13.15062174086928
$ws.Range("L5").Value = 11.01425840595383
$ws.Range("M5").Value = 14.50878228266646
$ws.Range("O5").Value = 19.34790548338527
# Row 6
$ws.Range("B6").Value = 13.48014629143442
$ws.Range("C6").Value = 14.27261636913073
$ws.Range("D6").Value = 4.332510899375257
$ws.Range("F6").Value = 21.46157375860825
$ws.Range("G6").Value = 23.87183068261236
$ws.Range("H6").Value = 13.15339060003992
$ws.Range("L6").Value = 11.01543655735379
$ws.Range("M6").Value = 14.5035361083456
$ws.Range("O6").Value = 19.35261065473207
# Row 7
$ws.Range("B7").Value = 13.67279476053798
$ws.Range("C7").Value = 14.30528324984317
$ws.Range("D7").Value = 4.354951203408718
$ws.Range("F7").Value = 21.4389265774679
$ws.Range("G7").Value = 23.83588635472699
$ws.Range("H7").Value = 13.13440682982523
$ws.Range("L7").Value = 11.00741820800097
$ws.Range("M7").Value = 14.54005993396585
$ws.Range("O7").Value = 19.32043865205361
# Row 8
$ws.Range("B8").Value = 14.49213312342483
$ws.Range("C8").Value = 14.45109953723313
$ws.Range("D8").Value = 4.451663408048644
$ws.Range("F8").Value = 21.35546513478113
$ws.Range("G8").Value = 23.70193242752932
$ws.Range("H8").Value = 13.05645695417358
$ws.Range("L8").Value = 10.97597768704027
$ws.Range("M8").Value = 14.70376125345651
$ws.Range("O8").Value = 19.19055476764776
# Row 9
$ws.Range("B9").Value = 15.97599246812173
$ws.Range("C9").Value = 14.74157425103071
$ws.Range("D9").Value = 4.632086420997502
$ws.Range("F9").Value = 21.25152033527904
$ws.Range("G9").Value = 23.52921905097303
$ws.Range("H9").Value = 12.92454827597528
$ws.Range("L9").Value = 10.92849684508727
$ws.Range("M9").Value = 15.03209462964858
$ws.Range("O9").Value = 18.97959437726111
# Row 10
$ws.Range("B10").Value = 16.98262302469348
$ws.Range("C10").Value = 14.95643650700325
$ws.Range("D10").Value = 4.75825128703968
$ws.Range("F10").Value = 21.21195074082173
$ws.Range("G10").Value = 23.45825104729318
$ws.Range("H10").Value = 12.84043256773546
$ws.Range("L10").Value = 10.90225422584935
$ws.Range("M10").Value = 15.27616670910737
$ws.Range("O10").Value = 18.8515468459676
# Row 11
$ws.Range("B11").Value = 17.42101265523124
$ws.Range("C11").Value = 15.05421826370699
$ws.Range("D11").Value = 4.814109621368378
$ws.Range("F11").Value = 21.2020033240056
$ws.Range("G11").Value = 23.43830079486467
$ws.Range("H11").Value = 12.80495087477693
$ws.Range("L11").Value = 10.89219143357415
$ws.Range("M11").Value = 15.38747347858482
$ws.Range("O11").Value = 18.79920425900258
# Row 12
$ws.Range("B12").Value = 17.58411680271461
$ws.Range("C12").Value = 15.09122794981791
$ws.Range("D12").Value = 4.835029392180454
$ws.Range("F12").Value = 21.19939802161041
$ws.Range("G12").Value = 23.43253227416837
$ws.Range("H12").Value = 12.79191564846684
$ws.Range("L12").Value = 10.88865039138701
$ws.Range("M12").Value = 15.42963330105317
$ws.Range("O12").Value = 18.78023766700039
# Row 13
$ws.Range("B13").Value = 17.54911995010651
$ws.Range("C13").Value = 15.08325843640791
$ws.Range("D13").Value = 4.830534470592763
$ws.Range("F13").Value = 21.19990741406563
$ws.Range("G13").Value = 23.43369500792378
$ws.Range("H13").Value = 12.79470517476153
$ws.Range("L13").Value = 10.88940103377571
$ws.Range("M13").Value = 15.42055343478746
$ws.Range("O13").Value = 18.78428438304974
# Row 14
$ws.Range("B14").Value = 17.43449003712825
$ws.Range("C14").Value = 15.0572635724355
$ws.Range("D14").Value = 4.815835430113868
$ws.Range("F14").Value = 21.20176568816976
$ws.Range("G14").Value = 23.43779035217497
$ws.Range("H14").Value = 12.80387041955179
$ws.Range("L14").Value = 10.89189470923264
$ws.Range("M14").Value = 15.39094195284247
$ws.Range("O14").Value = 18.79762671405183
# Row 15
$ws.Range("B15").Value = 17.36389507012609
$ws.Range("C15").Value = 15.04133790006251
$ws.Range("D15").Value = 4.806801215331721
$ws.Range("F15").Value = 21.2030552920982
$ws.Range("G15").Value = 23.44053182429956
$ws.Range("H15").Value = 12.80953663040731
$ws.Range("L15").Value = 10.89345725223353
$ws.Range("M15").Value = 15.3728045439377
$ws.Range("O15").Value = 18.80591069352579
# Row 16
$ws.Range("B16").Value = 16.95357165625647
$ws.Range("C16").Value = 14.95004499569094
$ws.Range("D16").Value = 4.754568900194009
$ws.Range("F16").Value = 21.21276320378194
$ws.Range("G16").Value = 23.45980415674254
$ws.Range("H16").Value = 12.84280744252628
$ws.Range("L16").Value = 10.90294957606562
$ws.Range("M16").Value = 15.26889560084218
$ws.Range("O16").Value = 18.85508684867186
# Row 17
$ws.Range("B17").Value = 16.6967758162263
$ws.Range("C17").Value = 14.8940329141793
$ws.Range("D17").Value = 4.72212432266775
$ws.Range("F17").Value = 21.22078413467264
$ws.Range("G17").Value = 23.47479570512759
$ws.Range("H17").Value = 12.86393130893725
$ws.Range("L17").Value = 10.90925300263856
$ws.Range("M17").Value = 15.20520076485092
$ws.Range("O17").Value = 18.88677131040249
# Row 18
$ws.Range("B18").Value = 16.5472399012339
$ws.Range("C18").Value = 14.86182141438532
$ws.Range("D18").Value = 4.703319273850651
$ws.Range("F18").Value = 21.22615540297813
$ws.Range("G18").Value = 23.48457870664582
$ws.Range("H18").Value = 12.87634309901226
$ws.Range("L18").Value = 10.91305507157438
$ws.Range("M18").Value = 15.1685930208656
$ws.Range("O18").Value = 18.90555106364011
# Row 19
$ws.Range("B19").Value = 16.49629783954213
$ws.Range("C19").Value = 14.85091677221109
$ws.Range("D19").Value = 4.696927885347845
$ws.Range("F19").Value = 21.22810405578979
$ws.Range("G19").Value = 23.48808991684674
$ws.Range("H19").Value = 12.88059047751364
$ws.Range("L19").Value = 10.91437270432498
$ws.Range("M19").Value = 15.15620395172276
$ws.Range("O19").Value = 18.91200485152608
# Row 20
$ws.Range("B20").Value = 16.72430270011761
$ws.Range("C20").Value = 14.89999513456758
$ws.Range("D20").Value = 4.725593072877668
$ws.Range("F20").Value = 21.21985182992863
$ws.Range("G20").Value = 23.47307964797385
$ws.Range("H20").Value = 12.86165552627315
$ws.Range("L20").Value = 10.90856372522611
$ws.Range("M20").Value = 15.2119785218031
$ws.Range("O20").Value = 18.8833408970825
# Row 21
$ws.Range("B21").Value = 17.46823913708468
$ws.Range("C21").Value = 15.0648995756765
$ws.Range("D21").Value = 4.820159297705676
$ws.Range("F21").Value = 21.20118832047323
$ws.Range("G21").Value = 23.43653888305844
$ws.Range("H21").Value = 12.80116747818271
$ws.Range("L21").Value = 10.89115494287815
$ws.Range("M21").Value = 15.39963951534959
$ws.Range("O21").Value = 18.79368452010938
# Row 22
$ws.Range("B22").Value = 17.93748084955352
$ws.Range("C22").Value = 15.17255750582746
$ws.Range("D22").Value = 4.880603183511375
$ws.Range("F22").Value = 21.1957620172745
$ws.Range("G22").Value = 23.42307306881672
$ws.Range("H22").Value = 12.76397226726849
$ws.Range("L22").Value = 10.88134816266319
$ws.Range("M22").Value = 15.52233376048807
$ws.Range("O22").Value = 18.74007130702717
# Row 23
$ws.Range("B23").Value = 17.68861773687623
$ws.Range("C23").Value = 15.11511694301634
$ws.Range("D23").Value = 4.848471386543288
$ws.Range("F23").Value = 21.19803766538405
$ws.Range("G23").Value = 23.42930337696423
$ws.Range("H23").Value = 12.7836099483501
$ws.Range("L23").Value = 10.88643854603946
$ws.Range("M23").Value = 15.45685490918886
$ws.Range("O23").Value = 18.76822816573546
# Row 24
$ws.Range("B24").Value = 16.71186370833373
$ws.Range("C24").Value = 14.89729964220689
$ws.Range("D24").Value = 4.724025324158154
$ws.Range("F24").Value = 21.22027095806581
$ws.Range("G24").Value = 23.47385185212836
$ws.Range("H24").Value = 12.86268357473326
$ws.Range("L24").Value = 10.90887479263196
$ws.Range("M24").Value = 15.20891426053707
$ws.Range("O24").Value = 18.88489003069503
# Row 25
$ws.Range("B25").Value = 15.58870228010518
$ws.Range("C25").Value = 14.66264321265982
$ws.Range("D25").Value = 4.584346748626676
$ws.Range("F25").Value = 21.27320239470431
$ws.Range("G25").Value = 23.56619251558729
$ws.Range("H25").Value = 12.95798758444864
$ws.Range("L25").Value = 10.92849684508727
$ws.Range("M25").Value = 15.03209462964858
$ws.Range("O25").Value = 19.03195288503845
